# Updated cryptos list on Sat Jul  1 19:55:17 UTC 2023 with GitHub Actions
#
# Refresh the "Price" (column D) and "Volume(1h)" (column E) figures for the
# coinranking.com snapshot on Sheet1. Column D values are plain display text
# (e.g. "30.611.53", "0.4742", "1.000") rather than numbers, so cells whose
# new value would otherwise be auto-parsed by Excel as a genuine number are
# assigned with a leading apostrophe to force a text entry, matching the
# original text-typed cells exactly. Column E values already contain the
# "%" sign plus padding spaces, so they round-trip as text unmodified.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.611.53"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "1.923.05"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'246.66"
$ws.Range("E5").Value = "  +2.51%  "
$ws.Range("D7").Value = "'0.4742"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "'0.2884"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("D9").Value = "'0.06834"
$ws.Range("E9").Value = "  +3.59%  "
$ws.Range("D10").Value = "'105.32"
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("D11").Value = "'18.35"
$ws.Range("E11").Value = "  -4.01%  "
$ws.Range("D12").Value = "1.928.20"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "'0.07696"
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").Value = "'5.338"
$ws.Range("E14").Value = "  +4.25%  "
$ws.Range("D15").Value = "'0.6676"
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("D16").Value = "'291.53"
$ws.Range("E16").Value = "  -3.78%  "
$ws.Range("D17").Value = "30.608.16"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").Value = "'0.000007611"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "'12.95"
$ws.Range("D21").Value = "'5.551"
$ws.Range("E21").Value = "  +5.80%  "
$ws.Range("D22").Value = "2.171.20"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'6.422"
$ws.Range("E24").Value = "  +1.83%  "
$ws.Range("D25").Value = "'9.458"
$ws.Range("E25").Value = "  +2.86%  "
$ws.Range("D26").Value = "'167.47"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").Value = "'21.04"
$ws.Range("E27").Value = "  +6.63%  "
$ws.Range("D28").Value = "'2.114"
$ws.Range("E28").Value = "  +5.44%  "
$ws.Range("D29").Value = "'0.1071"
$ws.Range("E29").Value = "  -4.69%  "
$ws.Range("E30").Value = "  +3.24%  "
$ws.Range("D31").Value = "'4.180"
$ws.Range("E31").Value = "  +1.82%  "
$ws.Range("D32").Value = "'4.056"
$ws.Range("E32").Value = "  +3.46%  "
$ws.Range("D33").Value = "'0.05033"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("D34").Value = "'0.7374"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").Value = "'1.144"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'0.02064"
$ws.Range("E36").Value = "  +5.95%  "
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("D38").Value = "'2.686"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").Value = "'2.054"
$ws.Range("E39").Value = "  +0.73%  "
$ws.Range("D40").Value = "'111.29"
$ws.Range("E40").Value = "  +3.68%  "
$ws.Range("D41").Value = "'0.8723"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "'0.4362"
$ws.Range("E42").Value = "  +5.82%  "
$ws.Range("D43").Value = "'5.918"
$ws.Range("E43").Value = "  +2.09%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "'67.91"
$ws.Range("E45").Value = "  -3.03%  "
$ws.Range("D46").Value = "'7.268"
$ws.Range("E46").Value = "  +0.92%  "
$ws.Range("D47").Value = "'9.303"
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("D48").Value = "'48.14"
$ws.Range("E48").Value = "  +15.04%  "
$ws.Range("E49").Value = "  +3.10%  "
$ws.Range("D50").Value = "'34.97"
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").Value = "'0.2503"
$ws.Range("E51").Value = "  +11.59%  "
